$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 4.245694160461426
$ws.Range("B1").Value = 3.370534181594849
$ws.Range("C1").Value = 2.969667673110962
$ws.Range("D1").Value = 3.709372282028198
$ws.Range("E1").Value = 5.226315975189209
